$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra comment rows (3-5), keeping only the header + one data row
$ws.Range("A3:D5").EntireRow.Delete()

# Replace the remaining data row (row 2) with the new record.
# Force text storage for the date-looking string so it isn't auto-converted
# to a date serial number.
$ws.Range("B2").NumberFormat = "@"

$ws.Range("A2").Value = "Apple iPhone XS Max"
$ws.Range("B2").Value = "02/20/2019"
$ws.Range("C2").Value = "https://www.gsmarena.com/apple_iphone_xs_max-reviews-9319p1.php"
$ws.Range("D2").Value = "I just love using this phone`n"

# Restore the default (un-styled) cell appearance on B2 now that the value
# is safely stored as text.
$ws.Range("B2").Style = "Normal"

# Writing the embedded newline above triggers an automatic row-height bump;
# put row 2 back to the default (no custom height) like the rest of the sheet.
$ws.Rows.Item(2).AutoFit()
